# resultsALL.xlsx: drop the TD / GP_N / GP_T columns (BN:BP) which shifts
# every OI_* column three places to the left, then append the new
# OI_T_EPD column (summary data) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the three leading columns; everything to the right (all the OI_*
# columns) slides left by three, landing the table on A1:DC3.
$ws.Range("BN1:BP1").EntireColumn.Delete()

# Copy the formatting of the new last existing header cell (DC1, bold +
# border, style index 1) onto the new DD1 header before writing its text,
# so the appended column matches the rest of the header row.
$ws.Range("DC1").Copy()
$ws.Range("DD1").PasteSpecial(-4122)

$ws.Range("DD1").Value = "OI_T_EPD"
$ws.Range("DD2").Value = 249.672
$ws.Range("DD3").Value = 340.717
